# RMA Complete Flow (Issue Credit) - SO TO RMA Receipt To Create Credit Memo
# Added test cases for Recurring billing. Maintenance of other test cases.
#
# The "RMA Details Maintenance Grid" sheet holds the most-recently-generated
# RMA numbers / line ids used to drive the next automation run. Each re-run
# of the QA automation appends a fresh block of generated RMA identifiers
# and points the maintenance grid rows at the newest block.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RMA Details Maintenance Grid")

# Row 2 -> first RMA line of the newest batch (RMA-DBXX)
$ws.Range("E2").Value = "RMA-DBXX-001"
$ws.Range("F2").Value = "RMA-DBXX-1-1"
$ws.Range("J2").Value = "a7s5f000000xKBgAAM"

# Row 3 -> second RMA line of the newest batch
$ws.Range("E3").Value = "RMA-DBXX-002"
$ws.Range("F3").Value = "RMA-DBXX-1-2"
$ws.Range("J3").Value = "a7s5f000000xKBhAAM"

# Row 4 -> third RMA line of the newest batch
$ws.Range("E4").Value = "RMA-DBXX-003"
$ws.Range("F4").Value = "RMA-DBXX-1-3"
$ws.Range("J4").Value = "a7s5f000000xKBiAAM"
